# "Generate Report for Handback"
# Refreshes the handoff/handback timestamps on the handback-status report:
#  - Overview!G3          "Latest HO Xliff Generate Date" for the
#                          39f37232-...-352c5a2e387e.md row
#  - zh-cn!H3 / zh-cn!K3   "Correspond Handoff Datetime" / "Correspond
#                          Handback DateTime" for the zh-cn xliff row
#  - de-de!H3 / de-de!K3   same, for the de-de xliff row

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-21 00:20:19"

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-21 00:20:08"
$wsZhCn.Range("K3").Value = "2016-10-21 00:20:59"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-10-21 00:20:19"
$wsDeDe.Range("K3").Value = "2016-10-21 00:21:18"
